# Scheduled-runner refresh of Golem Profits market data.
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H:N) for the
# affected leve rows across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
# Cells that go blank (no longer populated by the source feed) are cleared
# rather than zeroed, matching upstream behaviour.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 39
$ws.Range("H39").Value = 676.5
$ws.Range("I39").Value = 963
$ws.Range("K39").Value = 2889
$ws.Range("M39").Value = -2593
# row 51
$ws.Range("H51").Value = 15249.5
$ws.Range("J51").Value = 15249.5
$ws.Range("L51").Value = 15249.5
$ws.Range("N51").Value = -16217.5
# row 52
$ws.Range("H52").Value = 601.1667
$ws.Range("I52").Value = 621.4
$ws.Range("K52").Value = 1864.2
$ws.Range("M52").Value = -1704.2
# row 55
$ws.Range("H55").Value = 5077.6
$ws.Range("I55").Value = 5722
$ws.Range("J55").Value = 2500
$ws.Range("K55").Value = 5722
$ws.Range("L55").Value = 2500
$ws.Range("M55").Value = -5508
$ws.Range("N55").Value = -2928
# row 64
$ws.Range("H64").Value = 1999.5
$ws.Range("I64").Value = 1999.5
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 1999.5
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = -1751.5
$ws.Range("N64").Value = ""
# row 67
$ws.Range("H67").Value = 1999.5
$ws.Range("I67").Value = 1999.5
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 1999.5
$ws.Range("L67").Value = 0
$ws.Range("M67").Value = -1141.5
$ws.Range("N67").Value = ""
# row 138
$ws.Range("H138").Value = 2643.7778
$ws.Range("I138").Value = 796.4
$ws.Range("K138").Value = 2389.2
$ws.Range("M138").Value = 2750.8
# row 141
$ws.Range("H141").Value = 1723.3334
$ws.Range("I141").Value = 1723.3334
$ws.Range("K141").Value = 5170.0002
$ws.Range("M141").Value = 9.999799999999595
$ws = $wb.Worksheets.Item("ARM")
# row 45
$ws.Range("H45").Value = 2286.7778
$ws.Range("I45").Value = 2286.7778
$ws.Range("K45").Value = 2286.7778
$ws.Range("M45").Value = -1909.7778
# row 92
$ws.Range("H92").Value = 128163
$ws.Range("J92").Value = 128163
$ws.Range("L92").Value = 128163
$ws.Range("N92").Value = -133155
# row 97
$ws.Range("H97").Value = 41667204
$ws.Range("I97").Value = 41667204
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 41667204
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -41666708
$ws.Range("N97").Value = ""
# row 102
$ws.Range("I102").Value = 787465.0600000001
$ws.Range("K102").Value = 787465.0600000001
$ws.Range("M102").Value = -785843.0600000001
$ws = $wb.Worksheets.Item("BSM")
# row 54
$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("M54").Value = ""
$ws = $wb.Worksheets.Item("CRP")
# row 25
$ws.Range("H25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("N25").Value = ""
# row 31
$ws.Range("H31").Value = 3295.9443
$ws.Range("I31").Value = 1956.8462
$ws.Range("J31").Value = 6777.6
$ws.Range("K31").Value = 1956.8462
$ws.Range("L31").Value = 6777.6
$ws.Range("M31").Value = -1661.8462
$ws.Range("N31").Value = -7367.6
# row 34
$ws.Range("H34").Value = 3295.9443
$ws.Range("I34").Value = 1956.8462
$ws.Range("J34").Value = 6777.6
$ws.Range("K34").Value = 1956.8462
$ws.Range("L34").Value = 6777.6
$ws.Range("M34").Value = -1754.8462
$ws.Range("N34").Value = -7181.6
# row 134
$ws.Range("H134").Value = 1132.8889
$ws.Range("I134").Value = 1024.5
$ws.Range("K134").Value = 3073.5
$ws.Range("M134").Value = -538.5
$ws = $wb.Worksheets.Item("CUL")
# row 12
$ws.Range("H12").Value = 286.4
$ws.Range("I12").Value = 3.3333333
$ws.Range("J12").Value = 407.7143
$ws.Range("K12").Value = 9.999999900000001
$ws.Range("L12").Value = 1223.1429
$ws.Range("M12").Value = 163.0000001
$ws.Range("N12").Value = -1569.1429
# row 23
$ws.Range("H23").Value = 183.66667
$ws.Range("I23").Value = 183.66667
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 551.00001
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -316.00001
$ws.Range("N23").Value = ""
# row 51
$ws.Range("H51").Value = 6399.5
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 6399.5
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 19198.5
$ws.Range("M51").Value = ""
$ws.Range("N51").Value = -20118.5
# row 58
$ws.Range("H58").Value = 1005
$ws.Range("I58").Value = 1005
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 3015
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -2887
$ws.Range("N58").Value = ""
# row 60
$ws.Range("H60").Value = 176.83333
$ws.Range("I60").Value = 176.25
$ws.Range("J60").Value = 178
$ws.Range("K60").Value = 528.75
$ws.Range("L60").Value = 534
$ws.Range("M60").Value = -277.75
$ws.Range("N60").Value = -1036
# row 62
$ws.Range("H62").Value = 2495
$ws.Range("I62").Value = 2495
$ws.Range("K62").Value = 7485
$ws.Range("M62").Value = -6799
# row 63
$ws.Range("H63").Value = 812
$ws.Range("I63").Value = 812
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 2436
$ws.Range("L63").Value = 0
$ws.Range("N63").Value = ""
$ws.Range("M63").Value = -1687
# row 65
$ws.Range("H65").Value = 2495
$ws.Range("I65").Value = 2495
$ws.Range("K65").Value = 22455
$ws.Range("M65").Value = -19023
# row 66
$ws.Range("H66").Value = 812
$ws.Range("I66").Value = 812
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 7308
$ws.Range("L66").Value = 0
$ws.Range("N66").Value = ""
$ws.Range("M66").Value = -3564
$ws = $wb.Worksheets.Item("GSM")
# row 97
$ws.Range("H97").Value = 358.75
$ws.Range("I97").Value = 358.75
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 358.75
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = 137.25
$ws.Range("N97").Value = ""
$ws = $wb.Worksheets.Item("LTW")
# row 35
$ws.Range("H35").Value = 7655.091
$ws.Range("I35").Value = 1030.1428
$ws.Range("J35").Value = 19248.75
$ws.Range("K35").Value = 1030.1428
$ws.Range("L35").Value = 19248.75
$ws.Range("M35").Value = -694.1428000000001
$ws.Range("N35").Value = -19920.75
$ws = $wb.Worksheets.Item("WVR")
# row 107
$ws.Range("H107").Value = 633.1667
$ws.Range("J107").Value = 966.6667
$ws.Range("L107").Value = 2900.0001
$ws.Range("N107").Value = -6740.0001
# row 132
$ws.Range("H132").Value = 2544.4443
$ws.Range("I132").Value = 1580
$ws.Range("K132").Value = 4740
$ws.Range("M132").Value = -2210
# row 136
$ws.Range("H136").Value = 1653.375
$ws.Range("I136").Value = 1449.4
$ws.Range("K136").Value = 4348.200000000001
$ws.Range("M136").Value = -1798.200000000001
